$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAC")
$ws.Name = "PAC Characteristics"
